# Add a new "Sandbox header" paragraph at the end of the content
# placeholder's text on slide 7 (bold text, 30pt space-before, no bullet).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Append a new paragraph (CR = paragraph break) after the existing text.
$cr = [char]13
[void]$tr.InsertAfter($cr + "Sandbox header")

# Grab just the newly created last paragraph so formatting doesn't bleed
# into the other paragraphs in this text box.
$count = $tr.Paragraphs().Count
$newPara = $tr.Paragraphs($count, 1)

# Space before 30pt (spcPts val is in hundredths of a point -> 3000).
$newPara.ParagraphFormat.SpaceBefore = 30
# Keep paragraph bullet-free, re-asserted after SpaceBefore so the
# resulting <a:pPr> lists <a:spcBef> before <a:buNone> (matches native
# PowerPoint paragraph-property ordering).
$newPara.ParagraphFormat.Bullet.Type = 1
$newPara.ParagraphFormat.Bullet.Visible = $false

$newPara.Font.Bold = $true
